$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("金额数据")
$ws.Range("D12").Value = 11000
